$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "WARNING: replacement not found: $find"
    }
    return $ok
}

# 1. Title (appears twice - main heading + bold repeat near footer)
Replace-Text "Play Book of Shadows for Free – Review of Gameplay & Bonus Features" "Play Book of Shadows Free: Review of Bonus Features and Graphics"

# 2. "What we like" bullets
Replace-Text "Expanding symbols during free spins" "Wide variety of bonus features"
Replace-Text "Playable on desktop, mobile, and tablet devices" "Horror-themed graphics create a chilling atmosphere"
Replace-Text "Option to buy Free Spins feature" "Playability on different devices"

# 3. "What we don't like" bullets
Replace-Text "Limited bonus features compared to other slot games" "Limited number of paylines"
Replace-Text "High volatility may not be suitable for all players" "Option to buy the Free Spins feature may not appeal to all players"

# 4. Meta description (italic paragraph near the end)
Replace-Text "Read our review of Book of Shadows and play for free today. Learn about the gameplay mechanics and bonus features, including the Free Spins and Nolimit bonus." "Review of Book of Shadows online slot game. Play for free and enjoy bonus features and horror-themed graphics."
